$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet / reporting date from 2022-02-26 to 2022-02-27
$ws.Name = "Through 2022-02-27"

# Update header label for the "through" column
$ws.Range("I1").Value = "2022 (through 02-27)"

# Update February 2022 carjacking count (I3) and yearly total (I14)
$ws.Range("I3").Value = 136
$ws.Range("I14").Value = 295
